$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 5323.875
$ws.Range("J7").Value = 5748
$ws.Range("L7").Value = 5748
$ws.Range("N7").Value = -5972
# Row 14
$ws.Range("H14").Value = 5323.875
$ws.Range("J14").Value = 5748
$ws.Range("L14").Value = 5748
$ws.Range("N14").Value = -6130
# Row 18
$ws.Range("H18").Value = 394.94736
$ws.Range("I18").Value = 394.94736
$ws.Range("K18").Value = 394.94736
$ws.Range("M18").Value = -110.94736
# Row 19
$ws.Range("H19").Value = 1177.25
$ws.Range("J19").Value = 1146.2222
$ws.Range("L19").Value = 1146.2222
$ws.Range("N19").Value = -1496.2222
# Row 21
$ws.Range("H21").Value = 9384
$ws.Range("I21").Value = 9558
$ws.Range("J21").Value = 8949
$ws.Range("K21").Value = 9558
$ws.Range("L21").Value = 8949
$ws.Range("M21").Value = -9090
$ws.Range("N21").Value = -9885
# Row 23
$ws.Range("H23").Value = 9384
$ws.Range("I23").Value = 9558
$ws.Range("J23").Value = 8949
$ws.Range("K23").Value = 9558
$ws.Range("L23").Value = 8949
$ws.Range("M23").Value = -9324
$ws.Range("N23").Value = -9417
# Row 33
$ws.Range("H33").Value = 6614.0625
$ws.Range("I33").Value = 10383.2
$ws.Range("J33").Value = 332.16666
$ws.Range("K33").Value = 10383.2
$ws.Range("L33").Value = 332.16666
$ws.Range("M33").Value = -10154.2
$ws.Range("N33").Value = -790.16666
# Row 62
$ws.Range("H62").Value = 76809.47
$ws.Range("I62").Value = 101826.45
$ws.Range("J62").Value = 8012.75
$ws.Range("K62").Value = 101826.45
$ws.Range("L62").Value = 8012.75
$ws.Range("M62").Value = -101202.45
$ws.Range("N62").Value = -9260.75
# Row 65
$ws.Range("H65").Value = 76809.47
$ws.Range("I65").Value = 101826.45
$ws.Range("J65").Value = 8012.75
$ws.Range("K65").Value = 509132.25
$ws.Range("L65").Value = 40063.75
$ws.Range("M65").Value = -506012.25
$ws.Range("N65").Value = -46303.75
# Row 86
$ws.Range("H86").Value = 4963.0713
$ws.Range("I86").Value = 4367.6
$ws.Range("J86").Value = 6451.75
$ws.Range("K86").Value = 4367.6
$ws.Range("L86").Value = 6451.75
$ws.Range("M86").Value = -3244.6
$ws.Range("N86").Value = -8697.75
# Row 89
$ws.Range("H89").Value = 4963.0713
$ws.Range("I89").Value = 4367.6
$ws.Range("J89").Value = 6451.75
$ws.Range("K89").Value = 21838
$ws.Range("L89").Value = 32258.75
$ws.Range("M89").Value = -16222
$ws.Range("N89").Value = -43490.75
# Row 112
$ws.Range("H112").Value = 3582.9048
$ws.Range("J112").Value = 3985.1177
$ws.Range("L112").Value = 11955.3531
$ws.Range("N112").Value = -14171.3531
# Row 132
$ws.Range("H132").Value = 3480.1943
$ws.Range("I132").Value = 2952.0435
$ws.Range("J132").Value = 4414.615
$ws.Range("K132").Value = 8856.130500000001
$ws.Range("L132").Value = 13243.845
$ws.Range("M132").Value = -6326.130500000001
$ws.Range("N132").Value = -18303.845
# Row 133
$ws.Range("H133").Value = 78000
$ws.Range("J133").Value = 78000
$ws.Range("L133").Value = 78000
$ws.Range("N133").Value = -88120
# Row 137
$ws.Range("H137").Value = 10114.333
$ws.Range("J137").Value = 16280.954
$ws.Range("L137").Value = 48842.862
$ws.Range("N137").Value = -53942.862
# Row 138
$ws.Range("H138").Value = 2254.1765
$ws.Range("I138").Value = 2238.7693
$ws.Range("J138").Value = 2259.4473
$ws.Range("K138").Value = 6716.3079
$ws.Range("L138").Value = 6778.341899999999
$ws.Range("M138").Value = -1576.3079
$ws.Range("N138").Value = -17058.3419

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2705.889
$ws.Range("I32").Value = 1878.1666
$ws.Range("K32").Value = 1878.1666
$ws.Range("M32").Value = -1591.1666
# Row 61
$ws.Range("H61").Value = 10799.4
$ws.Range("I61").Value = 2999.3333
$ws.Range("K61").Value = 2999.3333
$ws.Range("M61").Value = -2787.3333
# Row 132
$ws.Range("H132").Value = 7956.069
$ws.Range("I132").Value = 6904.591
$ws.Range("J132").Value = 11260.714
$ws.Range("K132").Value = 20713.773
$ws.Range("L132").Value = 33782.142
$ws.Range("M132").Value = -18183.773
$ws.Range("N132").Value = -38842.142
# Row 136
$ws.Range("H136").Value = 10799.4
$ws.Range("I136").Value = 2999.3333
$ws.Range("K136").Value = 8997.999899999999
$ws.Range("M136").Value = -6447.999899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
# Row 107
$ws.Range("H107").Value = 940.8333
$ws.Range("I107").Value = 865.8889
$ws.Range("K107").Value = 865.8889
$ws.Range("M107").Value = 1054.1111
# Row 132
$ws.Range("H132").Value = 77378.96000000001
$ws.Range("I132").Value = 50000
$ws.Range("J132").Value = 79759.74000000001
$ws.Range("K132").Value = 50000
$ws.Range("L132").Value = 79759.74000000001
$ws.Range("M132").Value = -44940
$ws.Range("N132").Value = -89879.74000000001
# Row 134
$ws.Range("H134").Value = 12411.218
$ws.Range("I134").Value = 9733.823
$ws.Range("J134").Value = 19997.166
$ws.Range("K134").Value = 29201.469
$ws.Range("L134").Value = 59991.49800000001
$ws.Range("M134").Value = -26666.469
$ws.Range("N134").Value = -65061.49800000001

$ws = $wb.Worksheets.Item("CRP")
# Row 26
$ws.Range("H26").Value = 10010.5
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 10010.5
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 10010.5
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -10584.5
# Row 31
$ws.Range("H31").Value = 1716.7872
$ws.Range("I31").Value = 2039.6428
$ws.Range("J31").Value = 1579.8182
$ws.Range("K31").Value = 2039.6428
$ws.Range("L31").Value = 1579.8182
$ws.Range("M31").Value = -1744.6428
$ws.Range("N31").Value = -2169.8182
# Row 34
$ws.Range("H34").Value = 1716.7872
$ws.Range("I34").Value = 2039.6428
$ws.Range("J34").Value = 1579.8182
$ws.Range("K34").Value = 2039.6428
$ws.Range("L34").Value = 1579.8182
$ws.Range("M34").Value = -1837.6428
$ws.Range("N34").Value = -1983.8182
# Row 48
$ws.Range("H48").Value = 13583.714
$ws.Range("J48").Value = 23737
$ws.Range("L48").Value = 23737
$ws.Range("N48").Value = -24689
# Row 107
$ws.Range("H107").Value = 661.64703
$ws.Range("I107").Value = 399.5
$ws.Range("J107").Value = 1036.1428
$ws.Range("K107").Value = 399.5
$ws.Range("L107").Value = 1036.1428
$ws.Range("M107").Value = 1520.5
$ws.Range("N107").Value = -4876.1428
# Row 132
$ws.Range("H132").Value = 41772.28
$ws.Range("I132").Value = 35703.93
$ws.Range("K132").Value = 107111.79
$ws.Range("M132").Value = -104581.79

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 3906.375
$ws.Range("I3").Value = 3906.375
$ws.Range("K3").Value = 11719.125
$ws.Range("M3").Value = -11607.125
# Row 14
$ws.Range("H14").Value = 4574.16
$ws.Range("I14").Value = 4574.16
$ws.Range("K14").Value = 13722.48
$ws.Range("M14").Value = -13549.48
# Row 32
$ws.Range("H32").Value = 815.6667
$ws.Range("I32").Value = 773.5
$ws.Range("J32").Value = 900
$ws.Range("K32").Value = 2320.5
$ws.Range("L32").Value = 2700
$ws.Range("M32").Value = -2037.5
$ws.Range("N32").Value = -3266
# Row 80
$ws.Range("H80").Value = 24505.363
$ws.Range("J80").Value = 16524.428
$ws.Range("L80").Value = 49573.284
$ws.Range("N80").Value = -51445.284
# Row 83
$ws.Range("H83").Value = 24505.363
$ws.Range("J83").Value = 16524.428
$ws.Range("L83").Value = 148719.852
$ws.Range("N83").Value = -158079.852
# Row 107
$ws.Range("H107").Value = 1411.359
$ws.Range("I107").Value = 1023.3
$ws.Range("J107").Value = 1545.1724
$ws.Range("K107").Value = 3069.9
$ws.Range("L107").Value = 4635.5172
$ws.Range("M107").Value = -1149.9
$ws.Range("N107").Value = -8475.5172
# Row 109
$ws.Range("H109").Value = 5757.6
$ws.Range("I109").Value = 5207.5386
$ws.Range("K109").Value = 15622.6158
$ws.Range("M109").Value = -14582.6158
# Row 116
$ws.Range("H116").Value = 627.3333
$ws.Range("I116").Value = 627.3333
$ws.Range("K116").Value = 1881.9999
$ws.Range("M116").Value = 1560.0001
# Row 131
$ws.Range("H131").Value = 3514.4348
$ws.Range("I131").Value = 1193.3334
$ws.Range("K131").Value = 3580.0002
$ws.Range("M131").Value = 1459.9998
# Row 133
$ws.Range("H133").Value = 1167.1428
$ws.Range("I133").Value = 1390
$ws.Range("J133").Value = 1000
$ws.Range("K133").Value = 4170
$ws.Range("L133").Value = 3000
$ws.Range("M133").Value = 890
$ws.Range("N133").Value = -13120

$ws = $wb.Worksheets.Item("GSM")
# Row 29
$ws.Range("H29").Value = 3833
$ws.Range("I29").Value = 749.5
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 749.5
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = -459.5
$ws.Range("N29").Value = -10580

$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Range("H13").Value = 5608.2856
$ws.Range("J13").Value = 4043
$ws.Range("L13").Value = 4043
$ws.Range("N13").Value = -4323
# Row 20
$ws.Range("H20").Value = 34932.637
$ws.Range("J20").Value = 44068.617
$ws.Range("L20").Value = 44068.617
$ws.Range("N20").Value = -44520.617
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
# Row 136
$ws.Range("H136").Value = 15877468
$ws.Range("I136").Value = 22226580
$ws.Range("J136").Value = 4690
$ws.Range("K136").Value = 66679740
$ws.Range("L136").Value = 14070
$ws.Range("M136").Value = -66677190
$ws.Range("N136").Value = -19170

$ws = $wb.Worksheets.Item("WVR")
# Row 22
$ws.Range("H22").Value = 4378.5
$ws.Range("I22").Value = 1749.5
$ws.Range("J22").Value = 7007.5
$ws.Range("K22").Value = 1749.5
$ws.Range("L22").Value = 7007.5
$ws.Range("M22").Value = -1456.5
$ws.Range("N22").Value = -7593.5
# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
# Row 136
$ws.Range("H136").Value = 2919.5898
$ws.Range("I136").Value = 2696.6667
$ws.Range("J136").Value = 3662.6667
$ws.Range("K136").Value = 8090.000100000001
$ws.Range("L136").Value = 10988.0001
$ws.Range("M136").Value = -5540.000100000001
$ws.Range("N136").Value = -16088.0001

